# "Responded to more posts." - fill in actual time spent for a few more
# discussion-question response tasks on the week4 sheet, and update the
# running-total formula's cached value accordingly (Excel recalculates
# this automatically once the inputs change).

$wb = $excel.ActiveWorkbook

# Nudge the workbook window position (best-effort; mirrors the xWindow/
# yWindow change in the saved workbook view).
$win = $wb.Windows.Item(1)
$win.Left = -25440
$win.Top = 1600

$ws = $wb.Worksheets.Item("week4")

# "Actual time length to complete" column (C) for a handful of rows that
# were previously still blank.
$ws.Range("C12").Value = 0.020833333333333332   # DQ1 response 4
$ws.Range("C13").Value = 0.013888888888888888   # DQ1 response 5
$ws.Range("C17").Value = 0.013888888888888888   # DQ2 response 3
$ws.Range("C20").Value = 0.4583333333333333     # Hand-in assignment

# Move the active cell selection to C14, matching where editing left off.
$ws.Range("C14").Select()
